# Scheduled-runner update: refresh Universalis market-price derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 419.66666
$ws.Range("I11").Value = 419.66666
$ws.Range("K11").Value = 419.66666
$ws.Range("M11").Value = -279.66666
$ws.Range("H17").Value = 568.7027
$ws.Range("J17").Value = 568.7027
$ws.Range("L17").Value = 1706.1081
$ws.Range("N17").Value = -2042.1081
$ws.Range("H28").Value = 1645
$ws.Range("I28").Value = 959.05884
$ws.Range("K28").Value = 959.05884
$ws.Range("M28").Value = -474.05884
$ws.Range("H33").Value = 710.3
$ws.Range("I33").Value = 589.5714
$ws.Range("K33").Value = 589.5714
$ws.Range("M33").Value = -360.5714
$ws.Range("H43").Value = 1464
$ws.Range("I43").Value = 1499
$ws.Range("J43").Value = 1394
$ws.Range("K43").Value = 1499
$ws.Range("L43").Value = 1394
$ws.Range("M43").Value = -1430
$ws.Range("N43").Value = -1532
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 1000
$ws.Range("K49").Value = 3000
$ws.Range("M49").Value = -2864
$ws.Range("H51").Value = 9471.076999999999
$ws.Range("J51").Value = 9433.1
$ws.Range("L51").Value = 9433.1
$ws.Range("N51").Value = -10401.1
$ws.Range("H55").Value = 385
$ws.Range("J55").Value = 555.5
$ws.Range("L55").Value = 555.5
$ws.Range("N55").Value = -983.5
$ws.Range("H86").Value = 3187.0908
$ws.Range("I86").Value = 3040.4443
$ws.Range("J86").Value = 3847
$ws.Range("K86").Value = 3040.4443
$ws.Range("L86").Value = 3847
$ws.Range("M86").Value = -1917.4443
$ws.Range("N86").Value = -6093
$ws.Range("H89").Value = 3187.0908
$ws.Range("I89").Value = 3040.4443
$ws.Range("J89").Value = 3847
$ws.Range("K89").Value = 15202.2215
$ws.Range("L89").Value = 19235
$ws.Range("M89").Value = -9586.2215
$ws.Range("N89").Value = -30467
$ws.Range("H116").Value = 200
$ws.Range("I116").Value = 200
$ws.Range("K116").Value = 200
$ws.Range("M116").Value = 3242
$ws.Range("H132").Value = 1434142.9
$ws.Range("I132").Value = 6250
$ws.Range("K132").Value = 18750
$ws.Range("M132").Value = -16220
$ws.Range("H135").Value = 8799.691999999999
$ws.Range("I135").Value = 1199.75
$ws.Range("J135").Value = 99999
$ws.Range("K135").Value = 10797.75
$ws.Range("L135").Value = 899991
$ws.Range("M135").Value = -8262.75
$ws.Range("N135").Value = -905061

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 14983.333
$ws.Range("I33").Value = 14983.333
$ws.Range("K33").Value = 14983.333
$ws.Range("M33").Value = -14654.333
$ws.Range("H74").Value = 1716.2667
$ws.Range("I74").Value = 1395.4166
$ws.Range("K74").Value = 1395.4166
$ws.Range("M74").Value = -521.4166
$ws.Range("H77").Value = 1716.2667
$ws.Range("I77").Value = 1395.4166
$ws.Range("K77").Value = 6977.083000000001
$ws.Range("M77").Value = -2609.083000000001
$ws.Range("H132").Value = 2312.2
$ws.Range("I132").Value = 2346.889
$ws.Range("K132").Value = 7040.667
$ws.Range("M132").Value = -4510.667
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 499
$ws.Range("I22").Value = 499
$ws.Range("K22").Value = 499
$ws.Range("M22").Value = -326
$ws.Range("H99").Value = 4119.8
$ws.Range("I99").Value = 4175
$ws.Range("J99").Value = 3899
$ws.Range("K99").Value = 4175
$ws.Range("L99").Value = 3899
$ws.Range("M99").Value = -2677
$ws.Range("N99").Value = -6895

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2281.32
$ws.Range("I31").Value = 1786.6923
$ws.Range("J31").Value = 2817.1667
$ws.Range("K31").Value = 1786.6923
$ws.Range("L31").Value = 2817.1667
$ws.Range("M31").Value = -1491.6923
$ws.Range("N31").Value = -3407.1667
$ws.Range("H34").Value = 2281.32
$ws.Range("I34").Value = 1786.6923
$ws.Range("J34").Value = 2817.1667
$ws.Range("K34").Value = 1786.6923
$ws.Range("L34").Value = 2817.1667
$ws.Range("M34").Value = -1584.6923
$ws.Range("N34").Value = -3221.1667
$ws.Range("H41").Value = 15571.429
$ws.Range("H58").Value = 2203.818
$ws.Range("J58").Value = 2500
$ws.Range("L58").Value = 2500
$ws.Range("N58").Value = -2906
$ws.Range("H74").Value = 74119
$ws.Range("J74").Value = 74119
$ws.Range("L74").Value = 74119
$ws.Range("N74").Value = -75867
$ws.Range("H77").Value = 74119
$ws.Range("J77").Value = 74119
$ws.Range("L77").Value = 222357
$ws.Range("N77").Value = -231093
$ws.Range("H134").Value = 3491
$ws.Range("I134").Value = 3540.1
$ws.Range("K134").Value = 10620.3
$ws.Range("M134").Value = -8085.299999999999
$ws.Range("H136").Value = 2203.818
$ws.Range("J136").Value = 2500
$ws.Range("L136").Value = 7500
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 299.57895
$ws.Range("I12").Value = 391.66666
$ws.Range("J12").Value = 216.7
$ws.Range("K12").Value = 1174.99998
$ws.Range("L12").Value = 650.0999999999999
$ws.Range("M12").Value = -1001.99998
$ws.Range("N12").Value = -996.0999999999999
$ws.Range("H107").Value = 1072.8695
$ws.Range("J107").Value = 1304.9375
$ws.Range("L107").Value = 3914.8125
$ws.Range("N107").Value = -7754.8125
$ws.Range("H137").Value = 4249.75
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14809.5
$ws.Range("I80").Value = 5583.25
$ws.Range("J80").Value = 18500
$ws.Range("K80").Value = 5583.25
$ws.Range("L80").Value = 18500
$ws.Range("M80").Value = -4585.25
$ws.Range("N80").Value = -20496
$ws.Range("H83").Value = 14809.5
$ws.Range("I83").Value = 5583.25
$ws.Range("J83").Value = 18500
$ws.Range("K83").Value = 27916.25
$ws.Range("L83").Value = 92500
$ws.Range("M83").Value = -22924.25
$ws.Range("N83").Value = -102484
$ws.Range("H126").Value = 2917.5
$ws.Range("I126").Value = 2012
$ws.Range("J126").Value = 3098.6
$ws.Range("K126").Value = 6036
$ws.Range("L126").Value = 9295.799999999999
$ws.Range("M126").Value = -3566
$ws.Range("N126").Value = -14235.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3587.25
$ws.Range("I46").Value = 1199
$ws.Range("J46").Value = 3928.4285
$ws.Range("K46").Value = 1199
$ws.Range("L46").Value = 3928.4285
$ws.Range("M46").Value = -1011
$ws.Range("N46").Value = -4304.4285
$ws.Range("H61").Value = 1291.6111
$ws.Range("I61").Value = 1291.6111
$ws.Range("K61").Value = 1291.6111
$ws.Range("M61").Value = -1089.6111
$ws.Range("H100").Value = 2988.889
$ws.Range("I100").Value = 2900
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2900
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2359
$ws.Range("N100").Value = -4082
$ws.Range("H113").Value = 1291.6111
$ws.Range("I113").Value = 1291.6111
$ws.Range("K113").Value = 1291.6111
$ws.Range("M113").Value = 878.3888999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 43999.5
$ws.Range("J93").Value = 68000
$ws.Range("L93").Value = 68000
$ws.Range("N93").Value = -72992
$ws.Range("H122").Value = 4749.354
$ws.Range("I122").Value = 5054.353
$ws.Range("J122").Value = 4008.6428
$ws.Range("K122").Value = 15163.059
$ws.Range("L122").Value = 12025.9284
$ws.Range("M122").Value = -12713.059
$ws.Range("N122").Value = -16925.9284
$ws.Range("H132").Value = 9256.5
$ws.Range("I132").Value = 9173.333000000001
$ws.Range("K132").Value = 27519.999
$ws.Range("M132").Value = -24989.999
